$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relabel the Stat/Category column (T) for rows 2-10 with clarified, ---
# --- capitalized and unit-suffixed labels (shared-string edits in the diff). ---
$ws.Range("T2").Value = "Families' lights on "
$ws.Range("T3").Value = "Rental and utility assistance (K)"
$ws.Range("T4").Value = "Didn’t not qualify for government support"
$ws.Range("T5").Value = "Client records "
$ws.Range("T6").Value = "Residents "
$ws.Range("T7").Value = "70% of below AMI (%)"
$ws.Range("T8").Value = "70% or below median income (%)"
$ws.Range("T9").Value = "Hispanic (%)"
$ws.Range("T10").Value = "Spanish Speakers (%)"

# --- Update the Stat values (S) that changed alongside the relabeling ---
$ws.Range("S3").Value = 900     # was 900000 (now expressed in thousands, see "(K)" label)
$ws.Range("S8").Value = 95      # rotated value for "70% or below median income (%)"
$ws.Range("S9").Value = 79.8    # rotated value for "Hispanic (%)"
$ws.Range("S10").Value = 74.3   # rotated value for "Spanish Speakers (%)"

# --- Widen column T (20th column) to fit the new, longer labels ---
# The engine's stored <col width> = ColumnWidth + 0.8333, so 46.1666667 -> 47
$ws.Columns.Item(20).ColumnWidth = 46.1666667

# --- Scroll the view right so column O is left-most visible, then select T8 ---
# (mirrors the new ethnicity/TAY graph area added off-screen to the right)
$ws.Application.ActiveWindow.ScrollColumn = 15
$ws.Range("T8").Select()
